# Logged Week 15 and simulated Week 16
# Updates the "R" (road) row (row 3) target-depth totals on both the
# OFF and DEF sheets to reflect the newly logged Week 15 data plus the
# simulated Week 16 data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 204
$wsOff.Range("C3").Value = 124
$wsOff.Range("D3").Value = 43
$wsOff.Range("E3").Value = 21

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 148
$wsDef.Range("C3").Value = 107
$wsDef.Range("D3").Value = 55
$wsDef.Range("E3").Value = 28
$wsDef.Range("F3").Value = 2
$wsDef.Range("G3").Value = 1
